# CS133JS Lab07 Instructions - groupC: minor updates and clarifications.
#
# Applies, via Word COM-interop calls, the same textual edits captured in
# the target unified diff:
#   1. "Upload the following to the " -> three runs:
#        "Upload the following " / "6 files " / "to the "
#   2. "The" + " " + "html " + "file for" (4 runs) -> single run
#        "The html file for"
#   3. "A zip file containing the four files" -> "The four files"
#   4. "A" / " code review " / "of your own code." -> "A code review of
#      your code " / "with the " / a "Prod" column sentence + new "." run
#   5. A new, empty trailing paragraph (holding the "_GoBack" bookmark) is
#      appended at the end of the document.

$d = $word.ActiveDocument

function Insert-RunXml($range, [string]$innerXml) {
    # Wrap a fragment of <w:p>...</w:p> markup in the minimal WordOpenXML
    # "single part" package Range.InsertXML expects, then insert it at the
    # (collapsed) range. When the range sits exactly at the start of an
    # existing paragraph, Word merges the inserted runs into that
    # paragraph instead of splitting it into a new one.
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="256">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        $innerXml + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1) "Upload the following to the " -> 3 runs with identical formatting.
# ---------------------------------------------------------------------
$p42 = $d.Paragraphs.Item(42)
$rng = $p42.Range
$rng.Find.Execute("Upload the following to the ")
$rng.Delete()
$p42 = $d.Paragraphs.Item(42)
$insPt = $d.Range($p42.Range.Start, $p42.Range.Start)
$xml1 = '<w:p>' + `
    '<w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Upload the following </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">6 files </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">to the </w:t></w:r>' + `
    '</w:p>'
Insert-RunXml $insPt $xml1

# ---------------------------------------------------------------------
# 2) "The" / " " / "html " / "file for" -> single run "The html file for"
# ---------------------------------------------------------------------
$p43 = $d.Paragraphs.Item(43)
$rng = $p43.Range
$rng.Find.Execute("The html file for")
$rng.Delete()
$p43 = $d.Paragraphs.Item(43)
$insPt = $d.Range($p43.Range.Start, $p43.Range.Start)
$xml2 = '<w:p><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>The html file for</w:t></w:r></w:p>'
Insert-RunXml $insPt $xml2

# ---------------------------------------------------------------------
# 3) "A zip file containing the four files" -> "The four files"
# ---------------------------------------------------------------------
$p44 = $d.Paragraphs.Item(44)
$rng = $p44.Range
$rng.Find.Execute("A zip file containing the four files")
$rng.Delete()
$p44 = $d.Paragraphs.Item(44)
$insPt = $d.Range($p44.Range.Start, $p44.Range.Start)
$xml3 = '<w:p><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>The four files</w:t></w:r></w:p>'
Insert-RunXml $insPt $xml3

# ---------------------------------------------------------------------
# 4) "A" / " code review " / "of your own code." (last, underlined run)
#    -> "A code review of your code " / "with the " /
#       a "Prod" column sentence (still underlined) + new "." run.
# ---------------------------------------------------------------------
$p45 = $d.Paragraphs.Item(45)
$textRng = $d.Range($p45.Range.Start, $p45.Range.End - 1)
$textRng.Delete()
$p45 = $d.Paragraphs.Item(45)
$insPt = $d.Range($p45.Range.Start, $p45.Range.Start)
$quoteOpen = [char]8220
$quoteClose = [char]8221
$xml4 = '<w:p>' + `
    '<w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">A code review of your code </w:t></w:r>' + `
    '<w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">with the </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/></w:rPr><w:t>' + $quoteOpen + 'Prod' + $quoteClose + ' column filled in by you</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>.</w:t></w:r>' + `
    '</w:p>'
Insert-RunXml $insPt $xml4

# ---------------------------------------------------------------------
# 5) Append a new, empty trailing paragraph carrying the "_GoBack"
#    bookmark (Word drops this at the last edit position on save).
# ---------------------------------------------------------------------
$p45 = $d.Paragraphs.Item(45)
$endPt = $d.Range($p45.Range.End, $p45.Range.End)
$xml5 = '<w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Insert-RunXml $endPt $xml5
